$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header for column E: new category "Errori individuazione configurazione"
$ws.Range("E1").Value = "Errori individuazione`n configurazione"

# Row 10 updates
$ws.Range("B10").Value = 8
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 19
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 19

# Row 11 updates
$ws.Range("B11").Value = 14
$ws.Range("C11").Value = 1
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 0

# Row 12 updates (was empty, now filled in)
$ws.Range("B12").Value = 7
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 250
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

# Row 13 updates (was empty, now filled in)
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 120
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Update selection to E2
$ws.Range("E2").Select()
